$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15
$ws.Range("H15").Value = 3563.3713
$ws.Range("I15").Value = 3563.3713
$ws.Range("K15").Value = 10690.1139
$ws.Range("M15").Value = -10521.1139

# Row 16
$ws.Range("H16").Value = 2551.9
$ws.Range("I16").Value = 2006.6666
$ws.Range("J16").Value = 2785.5715
$ws.Range("K16").Value = 2006.6666
$ws.Range("L16").Value = 2785.5715
$ws.Range("M16").Value = -1776.6666
$ws.Range("N16").Value = -3245.5715

# Row 39
$ws.Range("H39").Value = 402.82352
$ws.Range("I39").Value = 158.7
$ws.Range("J39").Value = 504.54166
$ws.Range("K39").Value = 476.1
$ws.Range("L39").Value = 1513.62498
$ws.Range("M39").Value = -180.1
$ws.Range("N39").Value = -2105.62498

# Row 53
$ws.Range("H53").Value = 612.53845
$ws.Range("I53").Value = 390.57144
$ws.Range("K53").Value = 390.57144
$ws.Range("M53").Value = 246.42856

# Row 101
$ws.Range("H101").Value = 2454.4443
$ws.Range("I101").Value = 2798.5715
$ws.Range("J101").Value = 1250
$ws.Range("K101").Value = 8395.7145
$ws.Range("L101").Value = 3750
$ws.Range("M101").Value = -6773.7145
$ws.Range("N101").Value = -6994

# Row 117
$ws.Range("H117").Value = 133000
$ws.Range("J117").Value = 133000
$ws.Range("L117").Value = 133000
$ws.Range("N117").Value = -142178

# Row 123
$ws.Range("H123").Value = 71995
$ws.Range("J123").Value = 71995
$ws.Range("L123").Value = 71995
$ws.Range("N123").Value = -81795

# Row 132
$ws.Range("H132").Value = 2257.2927
$ws.Range("I132").Value = 1977.3158
$ws.Range("J132").Value = 5803.6665
$ws.Range("K132").Value = 5931.9474
$ws.Range("L132").Value = 17410.9995
$ws.Range("M132").Value = -3401.9474
$ws.Range("N132").Value = -22470.9995

# Row 137
$ws.Range("H137").Value = 7228.875
$ws.Range("I137").Value = 1098
$ws.Range("K137").Value = 3294
$ws.Range("M137").Value = -744

$ws = $wb.Worksheets.Item("ARM")
# Row 25
$ws.Range("H25").Value = 10000
$ws.Range("J25").Value = 10000
$ws.Range("L25").Value = 10000
$ws.Range("N25").Value = -10804

# Row 30
$ws.Range("H30").Value = 3387.5
$ws.Range("I30").Value = 1183.3334
$ws.Range("K30").Value = 1183.3334
$ws.Range("M30").Value = -1033.3334

# Row 132
$ws.Range("H132").Value = 10899.462
$ws.Range("I132").Value = 6489.7144
$ws.Range("J132").Value = 16044.167
$ws.Range("K132").Value = 19469.1432
$ws.Range("L132").Value = 48132.501
$ws.Range("M132").Value = -16939.1432
$ws.Range("N132").Value = -53192.501

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 1813.4
$ws.Range("I105").Value = 1074.0834
$ws.Range("J105").Value = 2922.375
$ws.Range("K105").Value = 1074.0834
$ws.Range("L105").Value = 2922.375
$ws.Range("M105").Value = 672.9166
$ws.Range("N105").Value = -6416.375

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 2841.8333
$ws.Range("I7").Value = 69.77778000000001
$ws.Range("J7").Value = 11158
$ws.Range("K7").Value = 69.77778000000001
$ws.Range("L7").Value = 11158
$ws.Range("M7").Value = 43.22221999999999
$ws.Range("N7").Value = -11384

# Row 29
$ws.Range("H29").Value = 0
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

# Row 31
$ws.Range("H31").Value = 1677790.1
$ws.Range("I31").Value = 39102.75
$ws.Range("K31").Value = 39102.75
$ws.Range("M31").Value = -38807.75

# Row 34
$ws.Range("H34").Value = 1677790.1
$ws.Range("I34").Value = 39102.75
$ws.Range("K34").Value = 39102.75
$ws.Range("M34").Value = -38900.75

# Row 39
$ws.Range("H39").Value = 5000
$ws.Range("I39").Value = 5000
$ws.Range("K39").Value = 5000
$ws.Range("M39").Value = -4609

# Row 41
$ws.Range("H41").Value = 20490.727
$ws.Range("J41").Value = 61300
$ws.Range("L41").Value = 61300
$ws.Range("N41").Value = -62156

# Row 49
$ws.Range("H49").Value = 5000
$ws.Range("I49").Value = 5000
$ws.Range("K49").Value = 5000
$ws.Range("M49").Value = -4818

# Row 60
$ws.Range("H60").Value = 60798.332
$ws.Range("J60").Value = 70499.75
$ws.Range("L60").Value = 70499.75
$ws.Range("N60").Value = -71521.75

# Row 125
$ws.Range("H125").Value = 248002.75
$ws.Range("J125").Value = 248002.75
$ws.Range("L125").Value = 248002.75
$ws.Range("N125").Value = -252922.75

$ws = $wb.Worksheets.Item("CUL")
# Row 2
$ws.Range("H2").Value = 16764.541
$ws.Range("J2").Value = 28644.572
$ws.Range("L2").Value = 171867.432
$ws.Range("N2").Value = -172093.432

# Row 19
$ws.Range("H19").Value = 0
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()

# Row 36
$ws.Range("H36").Value = 15224.5
$ws.Range("I36").Value = 450
$ws.Range("J36").Value = 29999
$ws.Range("K36").Value = 1350
$ws.Range("L36").Value = 89997
$ws.Range("M36").Value = -1181
$ws.Range("N36").Value = -90335

# Row 40
$ws.Range("H40").Value = 382.42856
$ws.Range("I40").Value = 1046
$ws.Range("J40").Value = 117
$ws.Range("K40").Value = 4184
$ws.Range("L40").Value = 468
$ws.Range("M40").Value = -4115
$ws.Range("N40").Value = -606

# Row 46
$ws.Range("H46").Value = 1595.8334
$ws.Range("I46").Value = 800
$ws.Range("K46").Value = 2400
$ws.Range("M46").Value = -2309

# Row 47
$ws.Range("H47").Value = 27266.5
$ws.Range("I47").Value = 4533
$ws.Range("J47").Value = 50000
$ws.Range("K47").Value = 13599
$ws.Range("L47").Value = 150000
$ws.Range("M47").Value = -13168
$ws.Range("N47").Value = -150862

# Row 68
$ws.Range("H68").Value = 1402.3125
$ws.Range("I68").Value = 1162.1428
$ws.Range("K68").Value = 3486.4284
$ws.Range("M68").Value = -2675.4284

# Row 71
$ws.Range("H71").Value = 1402.3125
$ws.Range("I71").Value = 1162.1428
$ws.Range("K71").Value = 10459.2852
$ws.Range("M71").Value = -6403.2852

# Row 132
$ws.Range("H132").Value = 2608.8572
$ws.Range("I132").Value = 2602.4
$ws.Range("K132").Value = 23421.6
$ws.Range("M132").Value = -20891.6

$ws = $wb.Worksheets.Item("GSM")
# Row 3
$ws.Range("H3").Value = 2854.889
$ws.Range("J3").Value = 727.7143
$ws.Range("L3").Value = 727.7143
$ws.Range("N3").Value = -959.7143

# Row 80
$ws.Range("H80").Value = 3164
$ws.Range("I80").Value = 3192
$ws.Range("K80").Value = 3192
$ws.Range("M80").Value = -2194

# Row 83
$ws.Range("H83").Value = 3164
$ws.Range("I83").Value = 3192
$ws.Range("K83").Value = 15960
$ws.Range("M83").Value = -10968

# Row 132
$ws.Range("H132").Value = 111115700
$ws.Range("I132").Value = 111115700
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 333347100
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -333344570
$ws.Range("N132").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 46
$ws.Range("H46").Value = 3590.6667
$ws.Range("I46").Value = 3395.7058
$ws.Range("J46").Value = 4064.1428
$ws.Range("K46").Value = 3395.7058
$ws.Range("L46").Value = 4064.1428
$ws.Range("M46").Value = -3207.7058
$ws.Range("N46").Value = -4440.1428

$ws = $wb.Worksheets.Item("WVR")
# Row 64
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()

# Row 67
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()

# Row 93
$ws.Range("H93").Value = 100581.336
$ws.Range("I93").Value = 90000
$ws.Range("K93").Value = 90000
$ws.Range("M93").Value = -87504

# Row 107
$ws.Range("H107").Value = 16130315
$ws.Range("I107").Value = 20834652
$ws.Range("K107").Value = 62503956
$ws.Range("M107").Value = -62502036
